{"js": "const replacements = [\n  [\"2024-07-08 Monday\", \"2024-07-09 Tuesday\"],\n  [\"20\\u00F78=2, 4\", \"74\\u00F77=10, 4\"],\n  [\"88\\u00F79=9, 7\", \"87\\u00F79=9, 6\"],\n  [\"93\\u00F75=18, 3\", \"88\\u00F72=44, 0\"],\n  [\"72\\u00F76=12, 0\", \"45\\u00F77=6, 3\"],\n  [\"35\\u00F76=5, 5\", \"38\\u00F73=12, 2\"],\n  [\"95\\u00F79=10, 5\", \"45\\u00F72=22, 1\"],\n  [\"79\\u00F72=39, 1\", \"43\\u00F73=14, 1\"],\n  [\"13\\u00F78=1, 5\", \"92\\u00F78=11, 4\"],\n  [\"34\\u00F74=8, 2\", \"15\\u00F73=5, 0\"],\n  [\"26\\u00F74=6, 2\", \"45\\u00F77=6, 3\"],\n  [\"56\\u00F75=11, 1\", \"88\\u00F73=29, 1\"],\n  [\"62\\u00F79=6, 8\", \"90\\u00F78=11, 2\"],\n  [\"34\\u00F77=4, 6\", \"23\\u00F78=2, 7\"],\n  [\"89\\u00F76=14, 5\", \"14\\u00F74=3, 2\"],\n  [\"80\\u00F77=11, 3\", \"19\\u00F73=6, 1\"],\n  [\"83\\u00F76=13, 5\", \"24\\u00F74=6, 0\"],\n  [\"80\\u00F72=40, 0\", \"71\\u00F73=23, 2\"],\n  [\"26\\u00F79=2, 8\", \"82\\u00F72=41, 0\"],\n  [\"80\\u00F76=13, 2\", \"58\\u00F75=11, 3\"],\n  [\"39\\u00F79=4, 3\", \"67\\u00F72=33, 1\"],\n  [\"65\\u00F76=10, 5\", \"90\\u00F77=12, 6\"],\n  [\"19\\u00F74=4, 3\", \"29\\u00F77=4, 1\"],\n  [\"12\\u00F75=2, 2\", \"43\\u00F76=7, 1\"],\n  [\"97\\u00F74=24, 1\", \"96\\u00F73=32, 0\"],\n  [\"30\\u00F76=5, 0\", \"32\\u00F79=3, 5\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-07-08 Monday\", \"2024-07-09 Tuesday\"),\n    @(\"20\u00f78=2, 4\", \"74\u00f77=10, 4\"),\n    @(\"88\u00f79=9, 7\", \"87\u00f79=9, 6\"),\n    @(\"93\u00f75=18, 3\", \"88\u00f72=44, 0\"),\n    @(\"72\u00f76=12, 0\", \"45\u00f77=6, 3\"),\n    @(\"35\u00f76=5, 5\", \"38\u00f73=12, 2\"),\n    @(\"95\u00f79=10, 5\", \"45\u00f72=22, 1\"),\n    @(\"79\u00f72=39, 1\", \"43\u00f73=14, 1\"),\n    @(\"13\u00f78=1, 5\", \"92\u00f78=11, 4\"),\n    @(\"34\u00f74=8, 2\", \"15\u00f73=5, 0\"),\n    @(\"26\u00f74=6, 2\", \"45\u00f77=6, 3\"),\n    @(\"56\u00f75=11, 1\", \"88\u00f73=29, 1\"),\n    @(\"62\u00f79=6, 8\", \"90\u00f78=11, 2\"),\n    @(\"34\u00f77=4, 6\", \"23\u00f78=2, 7\"),\n    @(\"89\u00f76=14, 5\", \"14\u00f74=3, 2\"),\n    @(\"80\u00f77=11, 3\", \"19\u00f73=6, 1\"),\n    @(\"83\u00f76=13, 5\", \"24\u00f74=6, 0\"),\n    @(\"80\u00f72=40, 0\", \"71\u00f73=23, 2\"),\n    @(\"26\u00f79=2, 8\", \"82\u00f72=41, 0\"),\n    @(\"80\u00f76=13, 2\", \"58\u00f75=11, 3\"),\n    @(\"39\u00f79=4, 3\", \"67\u00f72=33, 1\"),\n    @(\"65\u00f76=10, 5\", \"90\u00f77=12, 6\"),\n    @(\"19\u00f74=4, 3\", \"29\u00f77=4, 1\"),\n    @(\"12\u00f75=2, 2\", \"43\u00f76=7, 1\"),\n    @(\"97\u00f74=24, 1\", \"96\u00f73=32, 0\"),\n    @(\"30\u00f76=5, 0\", \"32\u00f79=3, 5\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
